# Whiteboarding Template.xlsx - add "Package" column support to the
# Example sheet, and make it the active/selected sheet.
#
# Commit message: "Included new line..again... Included support for
# Pacakge column. Will not prefix arguments with in/out/io if package
# type is library"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Insert a new first column ("Package") in front of the existing
# Workflow/Description/... columns.
$ws.Columns.Item(1).Insert() | Out-Null

# Header for the new column.
$ws.Range("A1").Value = "Package"

# New trailing row describing a "Library" package example (added before
# the Performer/Dispatcher rows below so the shared-string table keeps
# the same insertion order as the authored workbook).
$ws.Range("A5").Value = "Library"
$ws.Range("B5").Value = "App_Module2"
$ws.Range("C5").Value = "text"
$ws.Range("D5").Value = "text"
$ws.Range("E5").Value = "text"
$ws.Range("F5").Value = "in_Argument:string; out_Argument:int; io_Argument: dictionary"
$ws.Range("G5").Value = "text"
$ws.Range("H5").Value = "text"

# Tag the pre-existing example rows with their package type.
$ws.Range("A2").Value = "Performer"
$ws.Range("A3").Value = "Dispatcher"
$ws.Range("A4").Value = "Dispatcher"

# Make "Example" the active/selected sheet with C4 selected (matches the
# new tab selection + cursor position recorded in the workbook).
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null

# Restore the window from its minimized state.
$win = $excel.Windows.Item(1)
$win.WindowState = -4143
